$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.650777220726013
$ws.Range("B1").Value = 2.378071069717407
$ws.Range("C1").Value = 4.526685237884521
$ws.Range("D1").Value = 4.464462280273438
$ws.Range("E1").Value = 1.442973732948303
